$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 <- old row 27 content (except AC, which becomes empty)
$ws.Range("A25").Value = 112195278
$ws.Range("B25").Value = 8377
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 106545
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("Z25").Value = "11:00"
$ws.Range("AB25").Value = "11:00"
$ws.Range("AC25").Value = ""

# Row 26 <- old row 25 content (AC26 gains the comment that used to be on row 25)
$ws.Range("A26").Value = 112194720
$ws.Range("B26").Value = 56446
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 100049
$ws.Range("F26").Value = "Spillkråka"
$ws.Range("G26").Value = "Dryocopus martius"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("Z26").Value = "10:56"
$ws.Range("AB26").Value = "10:56"
$ws.Range("AC26").Value = "Hördes i omgivningen långa rop"

# Row 27 <- old row 26 content
$ws.Range("A27").Value = 112196967
$ws.Range("B27").Value = 43473
$ws.Range("E27").Value = 101735
$ws.Range("F27").Value = "Jättesvampmal"
$ws.Range("G27").Value = "Scardia boletella"
$ws.Range("H27").Value = "(Fabricius, 1794)"
$ws.Range("Z27").Value = "12:30"
$ws.Range("AB27").Value = "12:30"

# Row 28: new Id and updated start/end times
$ws.Range("A28").Value = 112196450
$ws.Range("Z28").Value = "12:23"
$ws.Range("AB28").Value = "12:23"
